$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (Ano 2025) with refreshed Bibi faturamento figures
$ws.Range("B9").Value = 3516351.29
$ws.Range("C9").Value = 550180.24
$ws.Range("D9").Value = 4066531.53
$ws.Range("E9").Value = 13.52947188386856
$ws.Range("F9").Value = 86.47052811613143
$ws.Range("G9").Value = -46.82782312914848
$ws.Range("H9").Value = -36.49955614265165
$ws.Range("I9").Value = 35313
$ws.Range("J9").Value = 1501
$ws.Range("K9").Value = 36814
$ws.Range("L9").Value = 25404
$ws.Range("M9").Value = 160.0744579593765
$ws.Range("N9").Value = 9.286092802111456
